# Automatische test-sync: 2025-08-03 15:13:50
# Append the new test-mail log entry (row 27) to the "Logs" sheet and
# refresh the "Overig" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$row = 27

$logs.Cells.Item($row, 1).Value = "Bel jij klant Jansen even?"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Testmail #19: Bel jij klant Jansen even?"
$logs.Cells.Item($row, 4).Value = "Overig"
$logs.Cells.Item($row, 5).Value = "Geachte heer/mevrouw,`r`nDank voor uw bericht. Om u verder te kunnen helpen, ontvang ik graag meer informatie. Kunt u aangeven waarom klant Jansen gebeld dient te worden en welk specifiek onderwerp dit betreft? Zo kan ik ervoor zorgen dat de juiste persoon of afdeling contact met hem opneemt.`r`nMet vriendelijke groet,`r`n[Naam]`r`nE-mailassistent van [Bedrijfsnaam]"
$logs.Cells.Item($row, 6).Value = "2025-08-03 15:13:02"
$logs.Cells.Item($row, 7).Value = "Ja"
$logs.Cells.Item($row, 8).Value = "Nee"
$logs.Cells.Item($row, 9).Value = "Ja"
$logs.Cells.Item($row, 10).Value = "Nee"

# Settle the row back to its default (non-custom) height after writing the
# multi-line "Antwoord" text, matching the other rows in the sheet.
$logs.Rows.Item($row).AutoFit() | Out-Null

# Extend the conditional-formatting ranges (one block per column) so they
# keep covering the data through the newly added row.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $fcs = $logs.Range($col + "2").FormatConditions
    $newRange = $logs.Range($col + "2:" + $col + $row)
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# The new row is categorised "Overig" - bump that count on the Dashboard.
$dashboard.Range("B2").Value = 9
